# Add simple Buttons for GUI
# Replace the static example values in B2/D2/B3/D3/D4/D5 with real
# dropdown (list) data validation on B2 (ContractType choices) and
# D2 (RequestDepartment choices), and clear the now-redundant example
# text that used to live in the shared strings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old placeholder example values, keep the cell formatting.
$ws.Range("B2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()

# The two validation-driven cells switch from the "Text" number format
# (used while they held literal example text) back to the sheet's
# default "General" format now that they host a dropdown instead.
$ws.Range("B2").NumberFormat = "General"
$ws.Range("D2").NumberFormat = "General"

# Add the dropdown "buttons" (list-type data validation) driving the GUI.
$ws.Range("B2").Validation.Add(3, 1, 1, '",Full time,Part time"')
$ws.Range("D2").Validation.Add(3, 1, 1, '",Administration,Services,Production,Financial"')

$wb.Save()
